$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Task "Valeur nulle et abérante (modif class et sensor)" is now done -> flag flips 0 -> 1
$ws.Range("E7").Value = 1

# Update the last active selection to reflect where the user ended up working
$ws.Range("I20").Select()
